$d = $word.ActiveDocument

$d.Content.Find.Execute("65×47=", $true, $false, $false, $false, $false, $true, 1, $false, "11×43=", 2) | Out-Null
$d.Content.Find.Execute("45×89=", $true, $false, $false, $false, $false, $true, 1, $false, "43×25=", 2) | Out-Null
$d.Content.Find.Execute("36×96=", $true, $false, $false, $false, $false, $true, 1, $false, "75×91=", 2) | Out-Null
$d.Content.Find.Execute("20×28=", $true, $false, $false, $false, $false, $true, 1, $false, "43×16=", 2) | Out-Null
$d.Content.Find.Execute("45×63=", $true, $false, $false, $false, $false, $true, 1, $false, "72×12=", 2) | Out-Null
$d.Content.Find.Execute("18×73=", $true, $false, $false, $false, $false, $true, 1, $false, "50×22=", 2) | Out-Null
$d.Content.Find.Execute("59×81=", $true, $false, $false, $false, $false, $true, 1, $false, "75×79=", 2) | Out-Null
$d.Content.Find.Execute("61×55=", $true, $false, $false, $false, $false, $true, 1, $false, "87×42=", 2) | Out-Null
$d.Content.Find.Execute("71×45=", $true, $false, $false, $false, $false, $true, 1, $false, "97×55=", 2) | Out-Null
$d.Content.Find.Execute("71×79=", $true, $false, $false, $false, $false, $true, 1, $false, "11×50=", 2) | Out-Null
$d.Content.Find.Execute("13×18=", $true, $false, $false, $false, $false, $true, 1, $false, "53×75=", 2) | Out-Null
$d.Content.Find.Execute("37×11=", $true, $false, $false, $false, $false, $true, 1, $false, "51×72=", 2) | Out-Null
$d.Content.Find.Execute("17×20=", $true, $false, $false, $false, $false, $true, 1, $false, "37×21=", 2) | Out-Null
$d.Content.Find.Execute("28×93=", $true, $false, $false, $false, $false, $true, 1, $false, "47×73=", 2) | Out-Null
$d.Content.Find.Execute("36×13=", $true, $false, $false, $false, $false, $true, 1, $false, "78×36=", 2) | Out-Null
$d.Content.Find.Execute("84×52=", $true, $false, $false, $false, $false, $true, 1, $false, "81×73=", 2) | Out-Null
$d.Content.Find.Execute("94×12=", $true, $false, $false, $false, $false, $true, 1, $false, "48×40=", 2) | Out-Null
$d.Content.Find.Execute("86×37=", $true, $false, $false, $false, $false, $true, 1, $false, "75×60=", 2) | Out-Null
$d.Content.Find.Execute("97×50=", $true, $false, $false, $false, $false, $true, 1, $false, "77×79=", 2) | Out-Null
$d.Content.Find.Execute("41×70=", $true, $false, $false, $false, $false, $true, 1, $false, "74×30=", 2) | Out-Null
$d.Content.Find.Execute("79×71=", $true, $false, $false, $false, $false, $true, 1, $false, "91×65=", 2) | Out-Null
$d.Content.Find.Execute("97×69=", $true, $false, $false, $false, $false, $true, 1, $false, "68×47=", 2) | Out-Null
$d.Content.Find.Execute("41×62=", $true, $false, $false, $false, $false, $true, 1, $false, "69×95=", 2) | Out-Null
$d.Content.Find.Execute("89×66=", $true, $false, $false, $false, $false, $true, 1, $false, "94×85=", 2) | Out-Null
$d.Content.Find.Execute("46×27=", $true, $false, $false, $false, $false, $true, 1, $false, "96×65=", 2) | Out-Null
